# Atualização semestre ... 2023.2
$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# 1) Slide 16 ("Desenvolvimento Web HTML, CSS, JS e PHP" cover slide):
#    add a new right-aligned "E-mail: ..." paragraph right after the
#    existing "Professor M.Sc. Heleno Cardoso" paragraph.
# -----------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
$shp16 = $slide16.Shapes.Item(5)
$tr16 = $shp16.TextFrame.TextRange

# Append a new paragraph (CR) carrying the whole e-mail line; it
# inherits the pPr (algn="r" hangingPunct="1") and rPr (sz=2000 b=1,
# bg1 solidFill) of the paragraph it is typed after.
$tr16.InsertAfter([char]13 + "E-mail: helenocardosofilho@gmail.com")

$emailPara = $tr16.Paragraphs(2, 1)

$runA = $emailPara.Characters(1, 8)    # "E-mail: "
$runB = $emailPara.Characters(9, 19)   # "helenocardosofilho@"
$runC = $emailPara.Characters(28, 5)   # "gmail"
$runD = $emailPara.Characters(33, 4)   # ".com"

# Touching Font.Bold (re-asserting the already-bold state) forces the
# engine to materialize each of these as its own <a:r> run instead of
# leaving the paragraph as a single run.
$runA.Font.Bold = $true
$runB.Font.Bold = $true
$runC.Font.Bold = $true
$runD.Font.Bold = $true

# -----------------------------------------------------------------
# 2) Slide 3 ("Apresentação Pessoal"): the bullet-list text box.
# -----------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$shp3 = $slide3.Shapes.Item(2)
$tr3 = $shp3.TextFrame.TextRange

# 2a) Paragraph 1 ("Analista de Sistemas; Lider SCRUM; Consultor e
#     Docente") is split across three runs; collapse it back into a
#     single run (keeping the first run's formatting) now that the
#     concatenated text is unchanged. Route through a temporary value
#     first since re-assigning the identical text is otherwise a
#     no-op for the engine's change detection.
$bioPara = $tr3.Paragraphs(1, 1)
$bioPara.Text = "__tmp__"
$bioPara2 = $tr3.Paragraphs(1, 1)
$bioPara2.Text = "Analista de Sistemas; Lider SCRUM; Consultor e Docente"

# 2b) Paragraph 6 ("Engenheiro Eletricista 7º Semestre - Área1
#     (Trancado)") becomes "Engenheiro Eletricista 8º/9º Semestre -
#     Área1 (Trancado)", split into three runs. Replace just the "7º "
#     substring with "8º/9º " so the run break lands after the
#     trailing space, matching the target run layout.
$semesterPara = $tr3.Paragraphs(6, 1)
$oldYear = $semesterPara.Characters(24, 3)   # "7º "
$oldYear.Text = "8º/9º "

Write-Host "edits applied"
